$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5: Sending cluster = ECs (values updated per Dr Hou advice)
# Rows 6-9: new rows, Sending cluster = M2

# Row 2
$ws.Range('A2').Value = 'ECs'
$ws.Range('B2').Value = 'S100a9'
$ws.Range('C2').Value = 'Tlr4'
$ws.Range('D2').Value = 'ECs'
$ws.Range('E2').Value = 1
$ws.Range('F2').Value = 0.3333333333333333
$ws.Range('G2').Value = 0.06313966666666666
$ws.Range('H2').Value = 0.189419
$ws.Range('I2').Value = 0.02041955056700611
$ws.Range('J2').Value = 0.02041955056700612
$ws.Range('K2').Value = 3
$ws.Range('L2').Value = 1
$ws.Range('M2').Value = 8.714516333333334
$ws.Range('N2').Value = 26.143549
$ws.Range('O2').Value = 0.1832255053237971
$ws.Range('P2').Value = 0.1832255053237971
$ws.Range('Q2').Value = 0.5502316564478889
$ws.Range('R2').Value = 4.952084908031
$ws.Range('S2').Value = 0.003741382471124522
$ws.Range('T2').Value = 0.003741382471124523

# Row 3
$ws.Range('A3').Value = 'ECs'
$ws.Range('B3').Value = 'S100a9'
$ws.Range('C3').Value = 'Tlr4'
$ws.Range('D3').Value = 'FAPs'
$ws.Range('E3').Value = 1
$ws.Range('F3').Value = 0.3333333333333333
$ws.Range('G3').Value = 0.06313966666666666
$ws.Range('H3').Value = 0.189419
$ws.Range('I3').Value = 0.02041955056700611
$ws.Range('J3').Value = 0.02041955056700612
$ws.Range('K3').Value = 3
$ws.Range('L3').Value = 1
$ws.Range('M3').Value = 16.48752133333333
$ws.Range('N3').Value = 49.462564
$ws.Range('O3').Value = 0.3466554324170239
$ws.Range('P3').Value = 0.346655432417024
$ws.Range('Q3').Value = 1.041016601146222
$ws.Range('R3').Value = 9.369149410316
$ws.Range('S3').Value = 0.007078548131566788
$ws.Range('T3').Value = 0.007078548131566793

# Row 4
$ws.Range('A4').Value = 'ECs'
$ws.Range('B4').Value = 'S100a9'
$ws.Range('C4').Value = 'Tlr4'
$ws.Range('D4').Value = 'M2'
$ws.Range('E4').Value = 1
$ws.Range('F4').Value = 0.3333333333333333
$ws.Range('G4').Value = 0.06313966666666666
$ws.Range('H4').Value = 0.189419
$ws.Range('I4').Value = 0.02041955056700611
$ws.Range('J4').Value = 0.02041955056700612
$ws.Range('K4').Value = 3
$ws.Range('L4').Value = 1
$ws.Range('M4').Value = 18.62376966666666
$ws.Range('N4').Value = 55.871309
$ws.Range('O4').Value = 0.391570739865005
$ws.Range('P4').Value = 0.391570739865005
$ws.Range('Q4').Value = 1.175898608830111
$ws.Range('R4').Value = 10.583087479471
$ws.Range('S4').Value = 0.007995698523233464
$ws.Range('T4').Value = 0.007995698523233466

# Row 5
$ws.Range('A5').Value = 'ECs'
$ws.Range('B5').Value = 'S100a9'
$ws.Range('C5').Value = 'Tlr4'
$ws.Range('D5').Value = 'sCs'
$ws.Range('E5').Value = 1
$ws.Range('F5').Value = 0.3333333333333333
$ws.Range('G5').Value = 0.06313966666666666
$ws.Range('H5').Value = 0.189419
$ws.Range('I5').Value = 0.02041955056700611
$ws.Range('J5').Value = 0.02041955056700612
$ws.Range('K5').Value = 3
$ws.Range('L5').Value = 1
$ws.Range('M5').Value = 3.735891666666667
$ws.Range('N5').Value = 11.207675
$ws.Range('O5').Value = 0.07854832239417409
$ws.Range('P5').Value = 0.0785483223941741
$ws.Range('Q5').Value = 0.2358829545361111
$ws.Range('R5').Value = 2.122946590825
$ws.Range('S5').Value = 0.001603921441081336
$ws.Range('T5').Value = 0.001603921441081337

# Row 6
$ws.Range('A6').Value = 'M2'
$ws.Range('B6').Value = 'S100a9'
$ws.Range('C6').Value = 'Tlr4'
$ws.Range('D6').Value = 'ECs'
$ws.Range('E6').Value = 3
$ws.Range('F6').Value = 1
$ws.Range('G6').Value = 3.028978666666667
$ws.Range('H6').Value = 9.086936
$ws.Range('I6').Value = 0.9795804494329938
$ws.Range('J6').Value = 0.9795804494329939
$ws.Range('K6').Value = 3
$ws.Range('L6').Value = 1
$ws.Range('M6').Value = 8.714516333333334
$ws.Range('N6').Value = 26.143549
$ws.Range('O6').Value = 0.1832255053237971
$ws.Range('P6').Value = 0.1832255053237971
$ws.Range('Q6').Value = 26.39608406398489
$ws.Range('R6').Value = 237.564756575864
$ws.Range('S6').Value = 0.1794841228526725
$ws.Range('T6').Value = 0.1794841228526726

# Row 7
$ws.Range('A7').Value = 'M2'
$ws.Range('B7').Value = 'S100a9'
$ws.Range('C7').Value = 'Tlr4'
$ws.Range('D7').Value = 'FAPs'
$ws.Range('E7').Value = 3
$ws.Range('F7').Value = 1
$ws.Range('G7').Value = 3.028978666666667
$ws.Range('H7').Value = 9.086936
$ws.Range('I7').Value = 0.9795804494329938
$ws.Range('J7').Value = 0.9795804494329939
$ws.Range('K7').Value = 3
$ws.Range('L7').Value = 1
$ws.Range('M7').Value = 16.48752133333333
$ws.Range('N7').Value = 49.462564
$ws.Range('O7').Value = 0.3466554324170239
$ws.Range('P7').Value = 0.346655432417024
$ws.Range('Q7').Value = 49.94035038487822
$ws.Range('R7').Value = 449.463153463904
$ws.Range('S7').Value = 0.3395768842854571
$ws.Range('T7').Value = 0.3395768842854572

# Row 8
$ws.Range('A8').Value = 'M2'
$ws.Range('B8').Value = 'S100a9'
$ws.Range('C8').Value = 'Tlr4'
$ws.Range('D8').Value = 'M2'
$ws.Range('E8').Value = 3
$ws.Range('F8').Value = 1
$ws.Range('G8').Value = 3.028978666666667
$ws.Range('H8').Value = 9.086936
$ws.Range('I8').Value = 0.9795804494329938
$ws.Range('J8').Value = 0.9795804494329939
$ws.Range('K8').Value = 3
$ws.Range('L8').Value = 1
$ws.Range('M8').Value = 18.62376966666666
$ws.Range('N8').Value = 55.871309
$ws.Range('O8').Value = 0.391570739865005
$ws.Range('P8').Value = 0.391570739865005
$ws.Range('Q8').Value = 56.4110010132471
$ws.Range('R8').Value = 507.6990091192239
$ws.Range('S8').Value = 0.3835750413417714
$ws.Range('T8').Value = 0.3835750413417715

# Row 9
$ws.Range('A9').Value = 'M2'
$ws.Range('B9').Value = 'S100a9'
$ws.Range('C9').Value = 'Tlr4'
$ws.Range('D9').Value = 'sCs'
$ws.Range('E9').Value = 3
$ws.Range('F9').Value = 1
$ws.Range('G9').Value = 3.028978666666667
$ws.Range('H9').Value = 9.086936
$ws.Range('I9').Value = 0.9795804494329938
$ws.Range('J9').Value = 0.9795804494329939
$ws.Range('K9').Value = 3
$ws.Range('L9').Value = 1
$ws.Range('M9').Value = 3.735891666666667
$ws.Range('N9').Value = 11.207675
$ws.Range('O9').Value = 0.07854832239417409
$ws.Range('P9').Value = 0.0785483223941741
$ws.Range('Q9').Value = 11.31593615931111
$ws.Range('R9').Value = 101.8434254338
$ws.Range('S9').Value = 0.07694440095309274
$ws.Range('T9').Value = 0.07694440095309277
